# Weekly update: insert two new price rows at the top of the Espinaca
# (spinach) price block for Mercado Mayorista Lo Valledor de Santiago,
# shifting all existing data rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 238; this pushes the
# existing rows 238-352 down to 240-354 and copies formatting from the
# row above (238 <- 237).
$ws.Rows("238:239").Insert()

# Populate the two newly inserted rows with this week's data.

# Row 238
$ws.Range("A238").Value2 = 6
$ws.Range("B238").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C238").Value2 = "Metropolitana"
$ws.Range("D238").Value2 = 44460
$ws.Range("E238").Value2 = 13
$ws.Range("F238").Value2 = 100112012
$ws.Range("G238").Value2 = "Espinaca"
$ws.Range("H238").Value2 = "Sin especificar"
$ws.Range("I238").Value2 = "Primera"
$ws.Range("J238").Value2 = 450
$ws.Range("K238").Value2 = 4500
$ws.Range("L238").Value2 = 5000
$ws.Range("M238").Value2 = 4711
$ws.Range("N238").Value2 = "$/cuna 10 kilos"
$ws.Range("O238").Value2 = "Provincia de Chacabuco"
$ws.Range("P238").Value2 = 471
$ws.Range("Q238").Value2 = 10
$ws.Range("R238").Value2 = "Hortaliza"

# Row 239
$ws.Range("A239").Value2 = 6
$ws.Range("B239").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C239").Value2 = "Metropolitana"
$ws.Range("D239").Value2 = 44460
$ws.Range("E239").Value2 = 13
$ws.Range("F239").Value2 = 100112012
$ws.Range("G239").Value2 = "Espinaca"
$ws.Range("H239").Value2 = "Sin especificar"
$ws.Range("I239").Value2 = "Primera"
$ws.Range("J239").Value2 = 420
$ws.Range("K239").Value2 = 4500
$ws.Range("L239").Value2 = 5000
$ws.Range("M239").Value2 = 4702
$ws.Range("N239").Value2 = "$/cuna 10 kilos"
$ws.Range("O239").Value2 = "Región Metropolitana"
$ws.Range("P239").Value2 = 470
$ws.Range("Q239").Value2 = 10
$ws.Range("R239").Value2 = "Hortaliza"
